# Add Transport-IDs column on Activity Executions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Transport IDS (falls mit Transport) (kommagetrennt)"

# Excel's ColumnWidth setter adds ~5/6 of a character as internal padding
# before it round-trips into the stored <col width="..."> value, so back
# that padding out here to land on the authored width of 54.
$ws.Columns.Item(9).ColumnWidth = 53.16666666666666

$ws.Range("I1").Select()
